$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 905.7059
$ws.Range("I137").Value = 812.25
$ws.Range("J137").Value = 1130
$ws.Range("K137").Value = 2436.75
$ws.Range("L137").Value = 3390
$ws.Range("M137").Value = 113.25
$ws.Range("N137").Value = -8490
$ws.Range("H138").Value = 4880.261
$ws.Range("I138").Value = 3370.3125
$ws.Range("J138").Value = 5215.8057
$ws.Range("K138").Value = 10110.9375
$ws.Range("L138").Value = 15647.4171
$ws.Range("M138").Value = -4970.9375
$ws.Range("N138").Value = -25927.4171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11227.528
$ws.Range("I32").Value = 9074.758
$ws.Range("K32").Value = 9074.758
$ws.Range("M32").Value = -8787.758
$ws.Range("H61").Value = 2194.2173
$ws.Range("I61").Value = 1951.5883
$ws.Range("J61").Value = 2881.6667
$ws.Range("K61").Value = 1951.5883
$ws.Range("L61").Value = 2881.6667
$ws.Range("M61").Value = -1739.5883
$ws.Range("N61").Value = -3305.6667
$ws.Range("H122").Value = 2275.0688
$ws.Range("I122").Value = 1513.9375
$ws.Range("J122").Value = 3211.8462
$ws.Range("K122").Value = 4541.8125
$ws.Range("L122").Value = 9635.5386
$ws.Range("M122").Value = -2091.8125
$ws.Range("N122").Value = -14535.5386
$ws.Range("H134").Value = 31400
$ws.Range("J134").Value = 37100
$ws.Range("L134").Value = 37100
$ws.Range("N134").Value = -47240
$ws.Range("H136").Value = 2194.2173
$ws.Range("I136").Value = 1951.5883
$ws.Range("J136").Value = 2881.6667
$ws.Range("K136").Value = 5854.7649
$ws.Range("L136").Value = 8645.0001
$ws.Range("M136").Value = -3304.7649
$ws.Range("N136").Value = -13745.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 80000000
$ws.Range("J19").Value = 80000000
$ws.Range("L19").Value = 80000000
$ws.Range("N19").Value = -80000346
$ws.Range("H107").Value = 949.3333
$ws.Range("I107").Value = 943
$ws.Range("K107").Value = 943
$ws.Range("M107").Value = 977

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H86").Value = 9618205
$ws.Range("I86").Value = 15627695
$ws.Range("J86").Value = 3019.5
$ws.Range("K86").Value = 15627695
$ws.Range("L86").Value = 3019.5
$ws.Range("M86").Value = -15626572
$ws.Range("N86").Value = -5265.5
$ws.Range("H89").Value = 9618205
$ws.Range("I89").Value = 15627695
$ws.Range("J89").Value = 3019.5
$ws.Range("K89").Value = 78138475
$ws.Range("L89").Value = 15097.5
$ws.Range("M89").Value = -78132859
$ws.Range("N89").Value = -26329.5
$ws.Range("H122").Value = 1379866
$ws.Range("I122").Value = 334300
$ws.Range("J122").Value = 2007205.6
$ws.Range("K122").Value = 1002900
$ws.Range("L122").Value = 6021616.8
$ws.Range("M122").Value = -1000450
$ws.Range("N122").Value = -6026516.8
$ws.Range("H132").Value = 864104.8
$ws.Range("I132").Value = 1389845.6
$ws.Range("J132").Value = 3801.6365
$ws.Range("K132").Value = 4169536.8
$ws.Range("L132").Value = 11404.9095
$ws.Range("M132").Value = -4167006.8
$ws.Range("N132").Value = -16464.9095
$ws.Range("H134").Value = 2669.72
$ws.Range("I134").Value = 2466.739
$ws.Range("J134").Value = 5004
$ws.Range("K134").Value = 7400.217
$ws.Range("L134").Value = 15012
$ws.Range("M134").Value = -4865.217
$ws.Range("N134").Value = -20082

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 836943.5
$ws.Range("I5").Value = 895.5
$ws.Range("J5").Value = 1351434.6
$ws.Range("K5").Value = 2686.5
$ws.Range("L5").Value = 4054303.8
$ws.Range("M5").Value = -2574.5
$ws.Range("N5").Value = -4054527.8
$ws.Range("H75").Value = 1280
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 1800
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 5400
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -7396
$ws.Range("H78").Value = 1280
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 1800
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 16200
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -26184
$ws.Range("H102").Value = 8265.5
$ws.Range("J102").Value = 8069.8887
$ws.Range("L102").Value = 24209.6661
$ws.Range("N102").Value = -29077.6661
$ws.Range("H123").Value = 1600
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H124").Value = 2420
$ws.Range("I124").Value = 500
$ws.Range("J124").Value = 2567.6924
$ws.Range("K124").Value = 1500
$ws.Range("L124").Value = 7703.0772
$ws.Range("M124").Value = 3410
$ws.Range("N124").Value = -17523.0772
$ws.Range("H135").Value = 836943.5
$ws.Range("I135").Value = 895.5
$ws.Range("J135").Value = 1351434.6
$ws.Range("K135").Value = 8059.5
$ws.Range("L135").Value = 12162911.4
$ws.Range("M135").Value = -5524.5
$ws.Range("N135").Value = -12167981.4
$ws.Range("H139").Value = 1422.2
$ws.Range("I139").Value = 980.35297
$ws.Range("K139").Value = 2941.05891
$ws.Range("M139").Value = 2198.94109

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3430.25
$ws.Range("I132").Value = 3182
$ws.Range("J132").Value = 3579.2
$ws.Range("K132").Value = 9546
$ws.Range("L132").Value = 10737.6
$ws.Range("M132").Value = -7016
$ws.Range("N132").Value = -15797.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 254.16667
$ws.Range("I55").Value = 141.83333
$ws.Range("J55").Value = 366.5
$ws.Range("K55").Value = 141.83333
$ws.Range("L55").Value = 366.5
$ws.Range("M55").Value = 31.16667000000001
$ws.Range("N55").Value = -712.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 33060.22
$ws.Range("I122").Value = 46883.137
$ws.Range("J122").Value = 2649.8
$ws.Range("K122").Value = 140649.411
$ws.Range("L122").Value = 7949.4
$ws.Range("M122").Value = -138199.411
$ws.Range("N122").Value = -12849.4
$ws.Range("H126").Value = 77599.53
$ws.Range("I126").Value = 87412.8
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 262238.4
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -259768.4
$ws.Range("N126").Value = -16940
$ws.Range("H136").Value = 4177.8843
$ws.Range("I136").Value = 2515.1082
$ws.Range("J136").Value = 6100.4688
$ws.Range("K136").Value = 7545.3246
$ws.Range("L136").Value = 18301.4064
$ws.Range("M136").Value = -4995.3246
$ws.Range("N136").Value = -23401.4064
